# Reprocess D2893, Edi06 and Edi09
# Updates H2O and N (nitrogen) values for samples D2893_2_bg (row 3) and
# D2893_3_bg_apf (row 4) across the mean, stdev and summary worksheets.

$wb = $excel.ActiveWorkbook

# --- "mean" sheet ---
$wsMean = $wb.Worksheets.Item("mean")
$wsMean.Range("J3").Value = 3.131   # H2O mean, D2893_2_bg
$wsMean.Range("K3").Value = 0.059   # N mean, D2893_2_bg
$wsMean.Range("J4").Value = 3.119   # H2O mean, D2893_3_bg_apf
$wsMean.Range("K4").Value = 0.071   # N mean, D2893_3_bg_apf

# --- "stdev" sheet ---
$wsStdev = $wb.Worksheets.Item("stdev")
$wsStdev.Range("J3").Value = 0.016  # H2O stdev, D2893_2_bg
$wsStdev.Range("K3").Value = 0.016  # N stdev, D2893_2_bg
$wsStdev.Range("J4").Value = 0.019  # H2O stdev, D2893_3_bg_apf
$wsStdev.Range("K4").Value = 0.019  # N stdev, D2893_3_bg_apf

# --- "summary" sheet ---
$wsSummary = $wb.Worksheets.Item("summary")
$wsSummary.Range("L3").Value = 3.131  # H2O mean, D2893_2_bg
$wsSummary.Range("M3").Value = 0.016  # H2O_sd, D2893_2_bg
$wsSummary.Range("T3").Value = 0.059  # N mean, D2893_2_bg
$wsSummary.Range("U3").Value = 0.016  # N_sd, D2893_2_bg

$wsSummary.Range("L4").Value = 3.119  # H2O mean, D2893_3_bg_apf
$wsSummary.Range("M4").Value = 0.019  # H2O_sd, D2893_3_bg_apf
$wsSummary.Range("T4").Value = 0.071  # N mean, D2893_3_bg_apf
$wsSummary.Range("U4").Value = 0.019  # N_sd, D2893_3_bg_apf
